# Ajout des formations, certif et langues dans l'export
#
# The three "recursive table" placeholder paragraphs for
# Competences Metier / Certifications / Langues are each immediately
# followed by an empty spacer paragraph. That spacer paragraph's
# paragraph-mark formatting (w:pPr/w:rPr) shrinks from 16pt (w:sz=32)
# down to 8pt (w:sz=16), matching the spacer that already follows the
# "Formations" table (which is already 8pt in this document).
#
# Because these paragraphs are empty (they contain only the paragraph
# mark, no run), Font.Size cannot be set on their Range directly in
# this host - the mutation silently does not stick. Instead we read
# back the exact WordOpenXML for that paragraph's Range (scoped to the
# /word/document.xml package part only, so styles.xml/numbering.xml
# etc. that are also echoed back in the package are left untouched),
# patch just the w:sz / w:szCs values for that one paragraph mark, and
# feed it back with Range.InsertXML so the paragraph mark keeps every
# other attribute (paraId, rsids, lang, bold, ...) unchanged.

$d = $word.ActiveDocument

function Set-ParaMarkHalfPointSize($para, [int]$oldHalfPoints, [int]$newHalfPoints) {
    $r = $para.Range
    $full = $r.WordOpenXML

    $partStartTag = '<pkg:part pkg:name="/word/document.xml"'
    $sIdx = $full.IndexOf($partStartTag)
    if ($sIdx -lt 0) {
        throw "could not locate /word/document.xml part in WordOpenXML"
    }
    $closeTag = "</pkg:part>"
    $eIdx = $full.IndexOf($closeTag, $sIdx) + $closeTag.Length

    $head = $full.Substring(0, $sIdx)
    $partXml = $full.Substring($sIdx, $eIdx - $sIdx)
    $tail = $full.Substring($eIdx)

    $oldSz = '<w:sz w:val="' + $oldHalfPoints + '"/>'
    $newSz = '<w:sz w:val="' + $newHalfPoints + '"/>'
    $oldSzCs = '<w:szCs w:val="' + $oldHalfPoints + '"/>'
    $newSzCs = '<w:szCs w:val="' + $newHalfPoints + '"/>'

    $newPartXml = $partXml.Replace($oldSz, $newSz)
    $newPartXml = $newPartXml.Replace($oldSzCs, $newSzCs)

    if ($newPartXml -eq $partXml) {
        throw "expected sz/szCs value not found on paragraph mark"
    }

    $newFull = $head + $newPartXml + $tail
    [void]$r.InsertXML($newFull)
}

function Get-ParagraphAfterUniqueText([string]$uniqueText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.MatchCase = $true
    $find.MatchWildcards = $false
    $find.Forward = $true
    $find.Wrap = 1
    $find.Text = $uniqueText
    $find.Replacement.Text = ""
    $ok = $find.Execute()
    if (-not $ok -or -not $find.Found) {
        throw "could not find text: $uniqueText"
    }
    $hitRange = $find.Parent
    $hitParagraph = $hitRange.Paragraphs(1)
    return $hitParagraph.Next()
}

$targets = @(
    "TABLEAU_RECURSIF_COMPETENCES_METIER",
    "TABLEAU_RECURSIF_CERTIFICATIONS",
    "TABLEAU_RECURSIF_LANGUES"
)

foreach ($placeholder in $targets) {
    $spacerPara = Get-ParagraphAfterUniqueText $placeholder
    Set-ParaMarkHalfPointSize $spacerPara 32 16
    Write-Output "Updated spacer paragraph after '$placeholder'"
}
